$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows (old rows 8-10, "MuSCs" as sending cluster set is fully
# replaced / the table shrinks from 3x3=9 data rows to 2x3=6 data rows)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pdpn"
$ws.Range("C2").Value = "Clec1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 45.835794
$ws.Range("H2").Value = 137.507382
$ws.Range("I2").Value = 0.8389444232146973
$ws.Range("J2").Value = 0.8389444232146973
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.075390666666667
$ws.Range("N2").Value = 3.226172
$ws.Range("O2").Value = 0.3593735692938006
$ws.Range("P2").Value = 0.3593735692938007
$ws.Range("Q2").Value = 49.29138506685599
$ws.Range("R2").Value = 443.6224656017041
$ws.Range("S2").Value = 0.3014944518097946
$ws.Range("T2").Value = 0.3014944518097947

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pdpn"
$ws.Range("C3").Value = "Clec1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 45.835794
$ws.Range("H3").Value = 137.507382
$ws.Range("I3").Value = 0.8389444232146973
$ws.Range("J3").Value = 0.8389444232146973
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.702433
$ws.Range("N3").Value = 5.107299
$ws.Range("O3").Value = 0.5689182942139039
$ws.Range("P3").Value = 0.5689182942139039
$ws.Range("Q3").Value = 78.032368286802
$ws.Range("R3").Value = 702.2913145812181
$ws.Range("S3").Value = 0.4772908301955731
$ws.Range("T3").Value = 0.4772908301955731

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pdpn"
$ws.Range("C4").Value = "Clec1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 45.835794
$ws.Range("H4").Value = 137.507382
$ws.Range("I4").Value = 0.8389444232146973
$ws.Range("J4").Value = 0.8389444232146973
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2145796666666667
$ws.Range("N4").Value = 0.6437390000000001
$ws.Range("O4").Value = 0.0717081364922955
$ws.Range("P4").Value = 0.0717081364922955
$ws.Range("Q4").Value = 9.835429397922
$ws.Range("R4").Value = 88.51886458129802
$ws.Range("S4").Value = 0.06015914120932964
$ws.Range("T4").Value = 0.06015914120932964

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pdpn"
$ws.Range("C5").Value = "Clec1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.799284
$ws.Range("H5").Value = 26.397852
$ws.Range("I5").Value = 0.1610555767853027
$ws.Range("J5").Value = 0.1610555767853027
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.075390666666667
$ws.Range("N5").Value = 3.226172
$ws.Range("O5").Value = 0.3593735692938006
$ws.Range("P5").Value = 0.3593735692938007
$ws.Range("Q5").Value = 9.462667886949333
$ws.Range("R5").Value = 85.164010982544
$ws.Range("S5").Value = 0.05787911748400599
$ws.Range("T5").Value = 0.057879117484006

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pdpn"
$ws.Range("C6").Value = "Clec1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.799284
$ws.Range("H6").Value = 26.397852
$ws.Range("I6").Value = 0.1610555767853027
$ws.Range("J6").Value = 0.1610555767853027
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.702433
$ws.Range("N6").Value = 5.107299
$ws.Range("O6").Value = 0.5689182942139039
$ws.Range("P6").Value = 0.5689182942139039
$ws.Range("Q6").Value = 14.980191457972
$ws.Range("R6").Value = 134.821723121748
$ws.Range("S6").Value = 0.09162746401833081
$ws.Range("T6").Value = 0.09162746401833081

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pdpn"
$ws.Range("C7").Value = "Clec1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.799284
$ws.Range("H7").Value = 26.397852
$ws.Range("I7").Value = 0.1610555767853027
$ws.Range("J7").Value = 0.1610555767853027
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2145796666666667
$ws.Range("N7").Value = 0.6437390000000001
$ws.Range("O7").Value = 0.0717081364922955
$ws.Range("P7").Value = 0.0717081364922955
$ws.Range("Q7").Value = 1.888147427625334
$ws.Range("R7").Value = 16.993326848628
$ws.Range("S7").Value = 0.01154899528296586
$ws.Range("T7").Value = 0.01154899528296586
